# Add a new "F" column (column 6) that duplicates column "E" (column 5,
# the first translated-language column) for each data row (rows 2-6),
# as part of adding the Norwegian translation column.
#
# Column F reuses the same cell formatting as E3:E6 (style index 3 in the
# original workbook - font "-apple-system", yellow-ish fill, left aligned,
# wrapped) for every new cell, including row 2 (where E2 itself happens to
# use a slightly different style). We copy that formatting via
# Copy/PasteSpecial so the existing style entry is reused instead of a new
# one being created, then copy each row's value (as a formula, to avoid
# re-typing literal text and to keep shared-string reuse) from E into F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 6
$srcCol = 5   # E
$dstCol = 6   # F

# Source of the cell formatting to stamp onto every new F cell: E3, which
# carries the "normal" translated-column style (cellXf index 3).
$formatSrc = $ws.Cells.Item(3, $srcCol)
$formatSrc.Copy()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dst = $ws.Cells.Item($r, $dstCol)
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, $srcCol)
    $dst = $ws.Cells.Item($r, $dstCol)
    $dst.Value = $src.Formula
}
